# Magento test data: append four new "Login Credentials" rows
# (rows 39-42) to match the refreshed ExtentReport test run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Credentials")

$rows = @(
    @("zijws@gmail.com", "LVhvl4596!", "Devika",       "Chattopadhyay D", "Devika Chattopadhyay D",       ","),
    @("pxbdg@gmail.com", "GPxg63717&", "Yogesh",       "Nambeesan D D",   "Yogesh Nambeesan D D",         ",,"),
    @("khgjd@gmail.com", "53upiu296%", "Abhirath",     "Devar",           "Abhirath Devar",               $null),
    @("uaiha@gmail.com", "US9wx3505#", "Brahmanandam", "Panicker D",      "Brahmanandam Panicker D",      ",")
)

$startRow = 39
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    if ($data[5] -ne $null) {
        $ws.Cells.Item($r, 6).Value = $data[5]
    }
}
